# Update Swap_Symbol_Details: avg_long (U) and avg_short (V) columns
# with refreshed print-out values per the commit "Updated Swap excel print outs".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 4;  U = -3.345;              V = -2.169 },
    @{ Row = 7;  U = 27.6091;             V = -53.17 },
    @{ Row = 8;  U = -3.5076;             V = -7.1422 },
    @{ Row = 9;  U = -78.70699999999999;  V = -157.271 },
    @{ Row = 10; U = 41.1387;             V = -66.9823 },
    @{ Row = 11; U = -2.1749;             V = -7.4385 },
    @{ Row = 12; U = 197.0945;            V = -284.81 },
    @{ Row = 13; U = -17.0854;            V = -24.5093 },
    @{ Row = 14; U = 0.06419999999999999; V = -6.3302 },
    @{ Row = 15; U = 0.1029;              V = -11.1824 },
    @{ Row = 17; U = -7.8636;             V = -13.3782 },
    @{ Row = 18; U = -14.3895;            V = -40.5616 },
    @{ Row = 19; U = -13.1448;            V = -50.679 },
    @{ Row = 20; U = 0.7593;              V = -22.991 },
    @{ Row = 21; U = -1.921;              V = 0.3755 },
    @{ Row = 22; U = -2.301;              V = -0.613 },
    @{ Row = 23; U = 0.3829;              V = -2.827 },
    @{ Row = 24; U = -1.112;              V = -1.366 },
    @{ Row = 25; U = -2.464;              V = -0.9350000000000001 },
    @{ Row = 26; U = -3.281;              V = -1.161 },
    @{ Row = 27; U = -2.755;              V = -1.666 },
    @{ Row = 28;                          V = -3.696 },
    @{ Row = 29; U = -0.362;              V = -2.248 },
    @{ Row = 30; U = 0.0459;              V = -1.686 },
    @{ Row = 31; U = -4.522;              V = 0.6693 },
    @{ Row = 32; U = -9.890000000000001;  V = -2.756 },
    @{ Row = 33; U = -5.429;              V = 0.1318 },
    @{ Row = 34; U = -6.135;              V = 0.9671 },
    @{ Row = 35; U = -0.875;              V = -2.567 },
    @{ Row = 37; U = -4.063;              V = -0.356 },
    @{ Row = 39; U = -4.217;              V = 0.06950000000000001 },
    @{ Row = 42; U = -6.749;              V = 0.8787 },
    @{ Row = 45; U = -7.885;              V = 0.2823 },
    @{ Row = 47; U = -5.5245;             V = 0.5821 },
    @{ Row = 48; U = -3.518;              V = -5.194 },
    @{ Row = 49; U = -4.938;              V = -3.355 },
    @{ Row = 50; U = 0.3764;              V = -6.898 },
    @{ Row = 52; U = -2.572;              V = -4.423 },
    @{ Row = 54; U = -5.7676;             V = -4.7051 },
    @{ Row = 56; U = -5.197;              V = -3.017 },
    @{ Row = 57; U = -3.881;              V = -3.026 },
    @{ Row = 62; U = -1.496;              V = -1.169 },
    @{ Row = 63; U = 0.6382;              V = -3.108 },
    @{ Row = 64; U = -0.466;              V = -1.806 },
    @{ Row = 65; U = -2.4;                V = -1.816 },
    @{ Row = 75; U = -1.4613;             V = -3.5601 }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey('U')) {
        $ws.Range("U$row").Value = $u.U
    }
    if ($u.ContainsKey('V')) {
        $ws.Range("V$row").Value = $u.V
    }
}
